$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap F:V between row 36 and row 37 ---
$tmp36 = $ws.Range("F36:V36").Value2
$tmp37 = $ws.Range("F37:V37").Value2
$ws.Range("F36:V36").Value2 = $tmp37
$ws.Range("F37:V37").Value2 = $tmp36

# --- Step 2: add 6 new rows (67-72), copying format from row 66 ---
$ws.Range("A66:V66").Copy($ws.Range("A67:V72"))

# Row 67
$ws.Range("A67").Value2 = 66
$ws.Range("B67").Value2 = "bosnia-and-herzegovina"
$ws.Range("C67").Value2 = "premijer-liga-bih"
$ws.Range("D67").Value2 = "2023-2024"
$ws.Range("E67").Value2 = 45226.58333333334
$ws.Range("F67").Value2 = "Igman K."
$ws.Range("G67").Value2 = 2
$ws.Range("H67").Value2 = "Zeljeznicar"
$ws.Range("I67").Value2 = 0
$ws.Range("J67").Value2 = 2.65
$ws.Range("K67").Value2 = "26/10/2023 02:12"
$ws.Range("L67").Value2 = 2.73
$ws.Range("M67").Value2 = "27/10/2023 13:54"
$ws.Range("N67").Value2 = 2.95
$ws.Range("O67").Value2 = "26/10/2023 02:12"
$ws.Range("P67").Value2 = 2.93
$ws.Range("Q67").Value2 = "27/10/2023 13:55"
$ws.Range("R67").Value2 = 2.51
$ws.Range("S67").Value2 = "26/10/2023 02:12"
$ws.Range("T67").Value2 = 2.75
$ws.Range("U67").Value2 = "27/10/2023 13:55"
$ws.Range("V67").Value2 = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/igman-konjic-zeljeznicar/h0Yakg9O/"

# Row 68
$ws.Range("A68").Value2 = 67
$ws.Range("B68").Value2 = "bosnia-and-herzegovina"
$ws.Range("C68").Value2 = "premijer-liga-bih"
$ws.Range("D68").Value2 = "2023-2024"
$ws.Range("E68").Value2 = 45227.5625
$ws.Range("F68").Value2 = "Sloga Doboj"
$ws.Range("G68").Value2 = 3
$ws.Range("H68").Value2 = "FK Sarajevo"
$ws.Range("I68").Value2 = 0
$ws.Range("J68").Value2 = 3.28
$ws.Range("K68").Value2 = "27/10/2023 01:43"
$ws.Range("L68").Value2 = 3.06
$ws.Range("M68").Value2 = "28/10/2023 13:22"
$ws.Range("N68").Value2 = 3.19
$ws.Range("O68").Value2 = "27/10/2023 01:43"
$ws.Range("P68").Value2 = 3.38
$ws.Range("Q68").Value2 = "28/10/2023 13:22"
$ws.Range("R68").Value2 = 2.02
$ws.Range("S68").Value2 = "27/10/2023 01:43"
$ws.Range("T68").Value2 = 2.23
$ws.Range("U68").Value2 = "28/10/2023 13:22"
$ws.Range("V68").Value2 = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/sloga-doboj-fk-sarajevo/63MsExnP/"

# Row 69
$ws.Range("A69").Value2 = 68
$ws.Range("B69").Value2 = "bosnia-and-herzegovina"
$ws.Range("C69").Value2 = "premijer-liga-bih"
$ws.Range("D69").Value2 = "2023-2024"
$ws.Range("E69").Value2 = 45227.66666666666
$ws.Range("F69").Value2 = "Borac Banja Luka"
$ws.Range("G69").Value2 = 1
$ws.Range("H69").Value2 = "Siroki Brijeg"
$ws.Range("I69").Value2 = 0
$ws.Range("J69").Value2 = 1.41
$ws.Range("K69").Value2 = "27/10/2023 04:12"
$ws.Range("L69").Value2 = 1.47
$ws.Range("M69").Value2 = "28/10/2023 15:56"
$ws.Range("N69").Value2 = 4.08
$ws.Range("O69").Value2 = "27/10/2023 04:12"
$ws.Range("P69").Value2 = 3.38
$ws.Range("Q69").Value2 = "28/10/2023 15:56"
$ws.Range("R69").Value2 = 6.16
$ws.Range("S69").Value2 = "27/10/2023 04:12"
$ws.Range("T69").Value2 = 10.66
$ws.Range("U69").Value2 = "28/10/2023 15:56"
$ws.Range("V69").Value2 = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/borac-banja-luka-siroki-brijeg/rPD99dft/"

# Row 70
$ws.Range("A70").Value2 = 69
$ws.Range("B70").Value2 = "bosnia-and-herzegovina"
$ws.Range("C70").Value2 = "premijer-liga-bih"
$ws.Range("D70").Value2 = "2023-2024"
$ws.Range("E70").Value2 = 45228.58333333334
$ws.Range("F70").Value2 = "Zvijezda 09"
$ws.Range("G70").Value2 = 0
$ws.Range("H70").Value2 = "Velez Mostar"
$ws.Range("I70").Value2 = 0
$ws.Range("J70").Value2 = 3.08
$ws.Range("K70").Value2 = "28/10/2023 03:12"
$ws.Range("L70").Value2 = 3.09
$ws.Range("M70").Value2 = "29/10/2023 13:55"
$ws.Range("N70").Value2 = 3.08
$ws.Range("O70").Value2 = "28/10/2023 03:12"
$ws.Range("P70").Value2 = 3.65
$ws.Range("Q70").Value2 = "29/10/2023 13:55"
$ws.Range("R70").Value2 = 2.15
$ws.Range("S70").Value2 = "28/10/2023 03:12"
$ws.Range("T70").Value2 = 2.11
$ws.Range("U70").Value2 = "29/10/2023 13:55"
$ws.Range("V70").Value2 = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/zvijezda-09-velez-mostar/C41M6fva/"

# Row 71
$ws.Range("A71").Value2 = 70
$ws.Range("B71").Value2 = "bosnia-and-herzegovina"
$ws.Range("C71").Value2 = "premijer-liga-bih"
$ws.Range("D71").Value2 = "2023-2024"
$ws.Range("E71").Value2 = 45228.77777777778
$ws.Range("F71").Value2 = "Posusje"
$ws.Range("G71").Value2 = 2
$ws.Range("H71").Value2 = "Zrinjski"
$ws.Range("I71").Value2 = 2
$ws.Range("J71").Value2 = 3.55
$ws.Range("K71").Value2 = "28/10/2023 07:42"
$ws.Range("L71").Value2 = 3.89
$ws.Range("M71").Value2 = "29/10/2023 18:39"
$ws.Range("N71").Value2 = 3.16
$ws.Range("O71").Value2 = "28/10/2023 07:42"
$ws.Range("P71").Value2 = 3.29
$ws.Range("Q71").Value2 = "29/10/2023 18:39"
$ws.Range("R71").Value2 = 1.93
$ws.Range("S71").Value2 = "28/10/2023 07:42"
$ws.Range("T71").Value2 = 1.96
$ws.Range("U71").Value2 = "29/10/2023 18:39"
$ws.Range("V71").Value2 = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/posusje-zrinjski/Qm3E8G9n/"

# Row 72
$ws.Range("A72").Value2 = 71
$ws.Range("B72").Value2 = "bosnia-and-herzegovina"
$ws.Range("C72").Value2 = "premijer-liga-bih"
$ws.Range("D72").Value2 = "2023-2024"
$ws.Range("E72").Value2 = 45229.58333333334
$ws.Range("F72").Value2 = "Tuzla City"
$ws.Range("G72").Value2 = 2
$ws.Range("H72").Value2 = "GOSK Gabela"
$ws.Range("I72").Value2 = 2
$ws.Range("J72").Value2 = 1.68
$ws.Range("K72").Value2 = "29/10/2023 02:12"
$ws.Range("L72").Value2 = 1.65
$ws.Range("M72").Value2 = "30/10/2023 13:34"
$ws.Range("N72").Value2 = 3.6
$ws.Range("O72").Value2 = "29/10/2023 02:12"
$ws.Range("P72").Value2 = 3.82
$ws.Range("Q72").Value2 = "30/10/2023 13:34"
$ws.Range("R72").Value2 = 4.12
$ws.Range("S72").Value2 = "29/10/2023 02:12"
$ws.Range("T72").Value2 = 4.94
$ws.Range("U72").Value2 = "30/10/2023 13:34"
$ws.Range("V72").Value2 = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/tuzla-city-nk-gosk-gabela/4v2I7zPh/"
